$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Actual output and Result for row 8 (test case 7)
$ws.Range("F8").Value = "It gets displayed as expected"
$ws.Range("G8").Value = "Pass"

# Update the selection / view state to match the edited cell
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("F9").Select()
